$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 687 entirely (the "「車のバッテリーが要る！」" post), causing all
# subsequent rows to shift up by one.
$ws.Rows.Item(687).Delete()
